# Insert a new data row at row 631 (pushing the existing rows 631..731 down
# to 632..732) and populate the new row with the latest weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 631:731 down to 632:732, carrying formatting from row 631.
$ws.Rows.Item(631).Insert()

$ws.Cells.Item(631, 1).Value = 3
$ws.Cells.Item(631, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(631, 3).Value = 'Coquimbo'
$ws.Cells.Item(631, 4).Value = 45218
$ws.Cells.Item(631, 5).Value = 5
$ws.Cells.Item(631, 6).Value = 100112017
$ws.Cells.Item(631, 7).Value = 'Apio'
$ws.Cells.Item(631, 8).Value = 'Americana (o)'
$ws.Cells.Item(631, 9).Value = 'Primera'
$ws.Cells.Item(631, 10).Value = 250
$ws.Cells.Item(631, 11).Value = 7500
$ws.Cells.Item(631, 12).Value = 8000
$ws.Cells.Item(631, 13).Value = 7760
$ws.Cells.Item(631, 14).Value = '$/docena de matas'
$ws.Cells.Item(631, 15).Value = 'Pan de Azúcar'
$ws.Cells.Item(631, 16).Value = 1293
$ws.Cells.Item(631, 17).Value = 6
$ws.Cells.Item(631, 18).Value = 'Hortaliza'
